# Applies the diff described:
#  - Summary sheet: update name, income, totals, ratio
#  - Assets sheet: rename vehicle description, update values
#  - Liabilities sheet: replace "Auto Loans / Vehicle Loan 1" row with
#    "Credit Cards / Credit Card Balance" values, remove the old
#    Credit Cards row, and move TOTAL LIABILITIES row up from row 4 to row 3

$wb = $excel.ActiveWorkbook

# ---------- Summary sheet ----------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Rana Al Qubaisi"
$summary.Range("B4").Value = 2244.15
$summary.Range("B6").Value = 61574
$summary.Range("B7").Value = 6976
$summary.Range("B8").Value = 54598
$summary.Range("B9").Value = 8.83

# ---------- Assets sheet ----------
$assets = $wb.Worksheets.Item("Assets")
$assets.Range("B2").Value = "Economy Car"
$assets.Range("C2").Value = 60339
$assets.Range("C3").Value = 1235
$assets.Range("C4").Value = 61574

# ---------- Liabilities sheet ----------
$liabilities = $wb.Worksheets.Item("Liabilities")

# Row 2 becomes the Credit Cards entry (replacing Auto Loans / Vehicle Loan 1)
$liabilities.Range("A2").Value = "Credit Cards"
$liabilities.Range("B2").Value = "Credit Card Balance"
$liabilities.Range("C2").Value = 6976
$liabilities.Range("D2").Value = 349
$liabilities.Range("E2").Value = 1

# Remove the old row 3 (previous Credit Cards row) entirely, shifting the
# TOTAL LIABILITIES row (previously row 4) up into row 3.
$liabilities.Rows.Item(3).Delete()

# Update the (now row 3) TOTAL LIABILITIES amount
$liabilities.Range("C3").Value = 6976
